$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new values in row 1 (A1 is a numeric-looking string, force text format)
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "1213213"
$ws.Range("B1").Value = "PuHJo0YDXT"
$ws.Range("C1").Value = "ztkinrkjpk"

# Remove the old row 2 (which previously held the data)
$ws.Rows.Item(2).Delete()
